$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new "Number of trials" column at D (shifts old D..K to E..L),
#    then move the data that used to live in K (now shifted to L) into D -
#    this reproduces "Number of trials" moving from the last column to right
#    after "Date Ended".
# ---------------------------------------------------------------------------
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("L1:L3").Cut($ws.Range("D1:D3"))

# ---------------------------------------------------------------------------
# 2. Append the 16 new per-interval response headers (L1:AA1), bold like the
#    rest of the header row.
# ---------------------------------------------------------------------------
$newHeaders = @(
    "PR - 3is", "OR - 3is", "TR - 3is", "SR - 3is",
    "PR - 6s",  "OR - 6s",  "TR - 6s",  "SR - 6s",
    "PR - 12s", "OR - 12s", "TR - 12s", "SR - 12s",
    "PR - 3fs", "OR - 3fs", "TR - 3fs", "SR - 3fs"
)
$col = 12
foreach ($h in $newHeaders) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}
$ws.Range("L1:AA1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Add the new "Test" session row (row 4).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Test"
$ws.Range("B4").Value = 44327.18941920334
$ws.Range("C4").Value = 44327.20417352678
$ws.Range("B4:C4").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# D4..AA4, in column order (Number of trials, Right.., Left.., Feeder.., Premature..,
# Omission.., Timed.., Perseverant.., then the 16 new interval columns)
$row4Values = @(0,0,0,120,0,120,0,0, 0,25,0,0, 0,34,0,0, 0,36,0,0, 0,25,0,0)
$col = 4
foreach ($v in $row4Values) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# 4. Selection / active cell moves to L1 after the edit.
# ---------------------------------------------------------------------------
$ws.Range("L1").Select() | Out-Null
